# confidence_rating.xlsx - label more than 50 points
#
# Sets column B ("rating") for a batch of previously-"low" rows to
# "moderate" or "high", adds a new "note" column (C) header, and annotates
# a handful of rows with "algae" / "human induced" notes.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New "note" column header
$ws.Range("C1").Value = "note"

# Rows whose rating moves from "low" to "moderate"
$moderateRows = @(22,31,40,42,80,86,131,168,171,182)
foreach ($r in $moderateRows) {
    $ws.Cells.Item($r, 2).Value = "moderate"
}

# Rows whose rating moves from "low" to "high"
$highRows = @(3,7,17,24,27,28,29,30,50,53,54,55,58,63,81,89,93,100,102,103,107,125,127,130,143,144,145,148,151,162,166,189,199,203,204)
foreach ($r in $highRows) {
    $ws.Cells.Item($r, 2).Value = "high"
}

# Rows annotated with a note in column C
$algaeRows = @(22,28,31,40,49)
foreach ($r in $algaeRows) {
    $ws.Cells.Item($r, 3).Value = "algae"
}

$humanRows = @(189)
foreach ($r in $humanRows) {
    $ws.Cells.Item($r, 3).Value = "human induced"
}

# Minor row-height tweaks picked up from the author's session
$ws.Rows.Item(129).RowHeight = 15.75
$ws.Rows.Item(204).RowHeight = 16.5

# Match the author's final selection/scroll position
$ws.Range("C104").Select()
